$d = $word.ActiveDocument

# 1. Replace the blank-line placeholder with the actual name.
$findRange = $d.Content
$findRange.Find.Execute("________", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Krizsa Mae", 2)

# After the replacement the paragraph text reads:
#   "Hello, my name is Krizsa Mae."
# "Krizsa Mae" spans characters 18-28, and the trailing "." is 28-29.
$nameStart = 18
$nameEnd   = 28

# 2. Move the (hidden) _GoBack bookmark so it sits right after the inserted
#    name and before the trailing period, splitting the run there.
$bm = $d.Bookmarks
$bm.Item("_GoBack").Delete()
$bm.Add("_GoBack", $d.Range($nameEnd, $nameEnd))

# 3. Also split "Hello, my name is " away from "Krizsa Mae" using a
#    throwaway bookmark as a run-splitting marker (added then immediately
#    removed), which leaves no extra formatting behind on the new runs.
$bm.Add("zzTempSplit", $d.Range($nameStart, $nameStart))
$bm.Item("zzTempSplit").Delete()
